$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.241.27"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "'3.515.27"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'610.66"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "'148.37"
$ws.Range("E6").Value = "  -2.16%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.480"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("D10").Value = "'0.142"
$ws.Range("E10").Value = "  -1.77%  "
$ws.Range("E11").Value = "  +6.21%  "
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").Value = "'4.110.19"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D16").Value = "'3.516.11"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "'67.288.08"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "'10.90"
$ws.Range("E19").Value = "  +8.97%  "
$ws.Range("D20").Value = "'6.36"
$ws.Range("E20").Value = "  -2.71%  "
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("D22").Value = "'436.81"
$ws.Range("E22").Value = "  -3.24%  "
$ws.Range("E23").Value = "  -3.08%  "
$ws.Range("D24").Value = "'80.14"
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("D25").Value = "'3.655.63"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -4.60%  "
$ws.Range("D28").Value = "'9.85"
$ws.Range("E28").Value = "  -1.66%  "
$ws.Range("D29").Value = "'8.30"
$ws.Range("E29").Value = "  -5.34%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  -5.42%  "
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("E33").Value = "  -2.69%  "
$ws.Range("D34").Value = "'25.60"
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("D35").Value = "'5.97"
$ws.Range("E35").Value = "  -4.26%  "
$ws.Range("D36").Value = "'1.82"
$ws.Range("E36").Value = "  -2.28%  "
$ws.Range("D37").Value = "'8.04"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("D40").Value = "'176.14"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("D41").Value = "'0.0902"
$ws.Range("E41").Value = "  -0.78%  "
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("E43").Value = "  -10.74%  "
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("D46").Value = "'28.40"
$ws.Range("E46").Value = "  -9.48%  "
$ws.Range("E47").Value = "  -5.90%  "
$ws.Range("E48").Value = "  -2.35%  "
$ws.Range("E49").Value = "  -3.22%  "
$ws.Range("D50").Value = "'0.993"
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("E51").Value = "  -2.69%  "
